$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.188.50'
$ws.Range('E2').Value = '  -5.72%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.557.74'
$ws.Range('E3').Value = '  -1.39%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '299.94'
$ws.Range('E5').Value = '  -2.59%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '94.56'
$ws.Range('E6').Value = '  -3.98%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.576'
$ws.Range('E7').Value = '  -2.59%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.552'
$ws.Range('E9').Value = '  -4.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.22'
$ws.Range('E10').Value = '  -6.01%  '
$ws.Range('E11').Value = '  -3.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.76'
$ws.Range('E12').Value = '  -3.59%  '
$ws.Range('E13').Value = '  +2.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.949.11'
$ws.Range('E14').Value = '  -1.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.527.35'
$ws.Range('E15').Value = '  -2.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.876'
$ws.Range('E16').Value = '  -3.57%  '
$ws.Range('E17').Value = '  -3.82%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.239.02'
$ws.Range('E18').Value = '  -5.84%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.14'
$ws.Range('E19').Value = '  +4.93%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0981'
$ws.Range('E20').Value = '  -2.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.63'
$ws.Range('E21').Value = '  -0.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '71.96'
$ws.Range('E22').Value = '  -1.41%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '260.18'
$ws.Range('E23').Value = '  -9.69%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.94'
$ws.Range('E24').Value = '  -2.81%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '29.65'
$ws.Range('E25').Value = '  +1.74%  '
$ws.Range('E26').Value = '  -4.89%  '
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.05'
$ws.Range('E28').Value = '  -5.68%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '37.68'
$ws.Range('E29').Value = '  -1.83%  '
$ws.Range('E30').Value = '  -2.86%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.00'
$ws.Range('E31').Value = '  -4.95%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '154.58'
$ws.Range('E32').Value = '  -2.27%  '
$ws.Range('E33').Value = '  -1.16%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.78'
$ws.Range('E34').Value = '  -1.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.39'
$ws.Range('E35').Value = '  -7.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0798'
$ws.Range('E36').Value = '  -4.40%  '
$ws.Range('E37').Value = '  -4.48%  '
$ws.Range('E38').Value = '  -2.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '17.00'
$ws.Range('E39').Value = '  +9.70%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '23.38'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.46'
$ws.Range('E41').Value = '  -1.76%  '
$ws.Range('E42').Value = '  -3.74%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.90'
$ws.Range('E43').Value = '  -2.27%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.087.58'
$ws.Range('E44').Value = '  -1.55%  '
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '85.93'
$ws.Range('E46').Value = '  -9.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.91'
$ws.Range('E47').Value = '  -3.47%  '
$ws.Range('E48').Value = '  +1.81%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.803.99'
$ws.Range('E49').Value = '  -1.38%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '104.90'
$ws.Range('E50').Value = '  -3.12%  '
$ws.Range('E51').Value = '  -1.67%  '
